# Update "想去人数" (F column) counts on the sheets that list event details.
# These values increased slightly from the previous crawl/generation.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, applied identically
# to both the "展览" and "全部类型" worksheets.
$updates = @{
    3  = 285
    4  = 293
    5  = 846
    6  = 14
    8  = 8362
    13 = 7
    15 = 24
    18 = 255
    19 = 720
    20 = 32
    21 = 82
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
